$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.055.02"
$ws.Range("D3").Value = "1.829.05"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.71"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6218"
$ws.Range("E6").Value = "  -6.29%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.42"
$ws.Range("E8").Value = "  +5.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07358"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2921"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.68"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07609"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("D13").Value = "1.830.66"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.961"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6618"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.10"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009063"
$ws.Range("E17").Value = "  +8.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.022"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").Value = "29.061.51"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "2.078.81"
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "225.29"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.36"
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.175"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.45"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.418"
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1357"
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.79"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.497"
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.054"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.035"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.200"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05237"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.843"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.151"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7332"
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.648"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("D39").Value = "1.287.02"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.748"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01782"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.323"
$ws.Range("E42").Value = "  +6.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8997"
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9994"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.77"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "1.976.57"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5114"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.95"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.707"
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3961"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.801"
$ws.Range("E51").Value = "  -0.10%  "
